$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mapping")

$ws.Range("A2").Value = -97.4535
$ws.Range("B2").Value = -97.3705

$ws.Range("A3").Value = 27.692
$ws.Range("B3").Value = 27.7659

$ws.Range("A4").Value = -96.5978
$ws.Range("B4").Value = -96.6795

$ws.Range("A5").Value = 28.4561
$ws.Range("B5").Value = 28.3829
